$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data change: officer David's age changed from 29 to 40 (cell C4)
$ws.Range("C4").Value = 40

# Leave the cursor where it was when the workbook was last saved
$ws.Range("E10").Select() | Out-Null
